# "Generate Report for Archive"
# Status moved on from handoff into translation: update the shared
# "Ready for handoff" status text to "In Translation" everywhere it is
# used (Overview zh-cn/de-de status columns, and each language sheet's
# own Status column), then let the narrower text's column shrink the
# columns that previously had to be wide enough for "Ready for handoff".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsZhCn.Range("C2").Value = "In Translation"
$wsDeDe.Range("C2").Value = "In Translation"

# The status columns narrow to fit the new (shorter) text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
